$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'knee pads snowboard'
    'youth under armour basketball tights'
    'black basketball rim'
    'mcdavid youth knee pads basketball'
    'nike compression pants men'
    'asics knee pad'
    'nike leggings basketball'
    'compression tights with knee pads'
    'combat knee pads'
    'uflex knee compression'
    'damascus knee pads'
    'mcdavid youth basketball knee pads'
    'mens pants black'
    'basketball pants girls'
    'compression capri leggings'
    'padded tights'
    'padded tights men'
    'crye knee pads'
    'dancing knee pads'
    'dodgeball knee pads'
    'uflex athletics knee compression sleeve'
    'gform knee pads'
    'nike tights for men'
    'insertable knee pads'
    'knee pads for scootering'
    'youth compression pants with pads'
    'red knee pads'
    'black volleyball knee pads'
    'men basketball pants'
    'compression tights with pads'
    'compression leggings with knee pads'
    'ski knee pads'
    'dead on knee pads'
    'black pads'
    'padded basketball compression pants'
    'mens padded leggings'
    'knee pads for skating'
    'padded compression tights'
    'basketball knee pads pants'
    'compression pants womens'
    'military knee pad'
    'drskin mens compression pants'
    'padded compression pants men basketball'
    'snowboard knee pad'
    'double knee pads'
    'knee pads light'
    'knee pad protection'
    'leggings with knee pads women'
    'padded tights men basketball'
    'mens leggings with knee pads'
    'knee pad under pants'
    'football pants adult xxl'
    'the best leggings'
    'sliding shorts mens'
    'snowboarding padded shorts men'
    'flexible knee pads'
    'basketball shorts for men pack'
    'leggings capri'
    'leggings spandex'
    'floor knee pads'
    'gym pads'
    'pants with knee pads kids'
    'basketball pants men'
    'volleyball knee pads nike black'
    'tights men'
    'hayabusa compression pants'
    'tough knee pads'
    'knee pad leggings'
    'goalkeeper pant'
    'knee compression men'
    'wrestling knee sleeves'
    'padded knee sleeves for basketball'
    'knee compression sleeve for squats'
    'youth leg sleeves for basketball'
    'tactical knee pads'
    'knee pads paintball'
    'bike knee pads'
    'knee pads biking'
    'goalie knee pads'
    'knee pads compression sleeve'
    'pant with knee pads'
    'long knee pads'
    'knee pads for girls'
    'mens pants with knee pads'
    'maroon knee pads'
    'pantalones con rodilleras'
    'pantalon con rodilleras'
    'knee padded pants men'
    'knee pad pants men'
    'baseball sliding pants mens'
    'men basketball knee pads'
    'baseball pants mens knee'
    'compression leggings basketball'
    'black mens baseball pants'
    'basketball pants for men'
    'work pants with knee pads'
    'knee pad for basketball youth'
    'basketball tights for men'
    'black athletic pants men'
    'youth basketball clothes'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

